$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with only Price (D) and Volume (E) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.411.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.917.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4815"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4068"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08220"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.013"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.30"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.946.12"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.247"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.47"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06822"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.63"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.006"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.434.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.658"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.74"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.174"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.155.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.02"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.113"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.44"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09597"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.673"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.548"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.374"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02282"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06106"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.181"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.84"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1845"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5590"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.954"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.427"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.31"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.37%  "

# --- Rows with only Volume (E) update ---
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.20%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.35%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.43%  "

# --- Rows that got re-ranked: Coin, Link, Price, Volume all change ---
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5987"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.12%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.067"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.23%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.280"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.56%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.402"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07595"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.71%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.40"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.97%  "
